$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Device name updated: S2400IB2 -> S2400IBH
$ws.Range("B4").Value = "NAME: S2400IBH"

# New curtain action added for CURTAIN_3 (OPEN -> WHOLE) on the S2400IBH remote link row
$ws.Range("C36").Value = "1: DEVICE CURTAIN_1 - CLOSE`n2: DEVICE CURTAIN_2 - OPEN`n3: DEVICE CURTAIN_3 - WHOLE`n4: SCENE Mixed Type`n5: GROUP DND"

# Leave the sheet scrolled/selected where the author ended up editing
$ws.Range("C35").Select()
